# Applies the diff described:
#  - rename sheet "fixed parameters" -> "static parameters"
#  - on "static parameters" sheet: add two new rows (operator/Robbert, device NO/TNO N77)
#  - widen column A on "static parameters" sheet and change selection
#  - on "dynamic parameters" sheet: clear cell B2 (was "C") and change selection to B2

$wb = $excel.ActiveWorkbook

$wsStatic = $wb.Worksheets.Item("fixed parameters")
$wsStatic.Name = "static parameters"

$wsStatic.Range("A6").Value = "operator"
$wsStatic.Range("B6").Value = "Robbert"
$wsStatic.Range("A7").Value = "device NO"
$wsStatic.Range("B7").Value = "TNO N77"

$wsStatic.Columns.Item(1).ColumnWidth = 20.75
$wsStatic.Range("A1:XFD1").Select() | Out-Null

$wsDynamic = $wb.Worksheets.Item("dynamic parameters")
$wsDynamic.Range("B2").ClearContents()
$wsDynamic.Range("B2").Select() | Out-Null
